$d = $word.ActiveDocument

# 1) Merge "Frontend Projekt mit HTML, CSS, " + "Javascript" (proofErr-split runs)
#    into a single run without the spell-check markers.
$d.Content.Find.Execute("Frontend Projekt mit HTML, CSS, Javascript", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Frontend Projekt mit HTML, CSS, Javascript", 2)

# 2) Merge "Unsere tägliche " + "Projektdoku" + " ist hier abgelegt und einsehbar: …"
#    into a single run without the spell-check markers.
$d.Content.Find.Execute("Unsere tägliche Projektdoku ist hier abgelegt und einsehbar: …", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Unsere tägliche Projektdoku ist hier abgelegt und einsehbar: …", 2)

# 3) Update the GitHub repository URL from Memory-Game to Memory-Spiel
$d.Content.Find.Execute("https://github.com/renewollny/Memory-Game", $true, $false, $false, $false, $false,
                         $true, 1, $false, "https://github.com/renewollny/Memory-Spiel", 2)
